$wb = $excel.ActiveWorkbook

# Cell value updates per sheet (scheduled market-data refresh).
# Keys are A1 cell refs; values are the new numeric values.

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$updates = @{
    "H9" = 899.6875
    "I9" = 903.2143
    "K9" = 903.2143
    "M9" = -734.2143
    "H17" = 2718.6
    "I17" = 990
    "J17" = 2910.6667
    "K17" = 2970
    "L17" = 8732.000100000001
    "M17" = -2802
    "N17" = -9068.000100000001
    "H18" = 718.25
    "I18" = 691
    "J18" = 800
    "K18" = 691
    "L18" = 800
    "M18" = -407
    "N18" = -1368
    "H86" = 900
    "J86" = 0
    "L86" = 0
    "H89" = 900
    "J89" = 0
    "L89" = 0
    "H137" = 2951.85
    "I137" = 1294.8334
    "K137" = 3884.5002
    "M137" = -1334.5002
    "H138" = 2093.6956
    "I138" = 1978.9048
    "K138" = 5936.7144
    "M138" = -796.7143999999998
    "H139" = 94000
    "I139" = 94000
    "J139" = 0
    "K139" = 94000
    "L139" = 0
    "N139" = -88860
    "H141" = 3584.9
    "I141" = 3584.9
    "K141" = 10754.7
    "M141" = -5574.700000000001
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$clears = @("N86", "N89", "M139")
foreach ($cellRef in $clears) {
    $ws.Range($cellRef).ClearContents()
}

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$updates = @{
    "H32" = 8780.941999999999
    "I32" = 8333.206
    "K32" = 8333.206
    "M32" = -8046.206
    "H45" = 2357.3
    "J45" = 0
    "L45" = 0
    "H61" = 5824.5
    "I61" = 5127.2964
    "K61" = 5127.2964
    "M61" = -4915.2964
    "H74" = 1693.75
    "I74" = 1766.6666
    "J74" = 600
    "K74" = 1766.6666
    "L74" = 600
    "M74" = -892.6666
    "N74" = -2348
    "H77" = 1693.75
    "I77" = 1766.6666
    "J77" = 600
    "K77" = 8833.333000000001
    "L77" = 3000
    "M77" = -4465.333000000001
    "N77" = -11736
    "H88" = 2878.1667
    "J88" = 2885.7273
    "L88" = 2885.7273
    "N88" = -3697.7273
    "H91" = 2878.1667
    "J91" = 2885.7273
    "L91" = 2885.7273
    "N91" = -5693.7273
    "H132" = 1375
    "I132" = 1373.5294
    "K132" = 4120.5882
    "M132" = -1590.5882
    "H136" = 5824.5
    "I136" = 5127.2964
    "K136" = 15381.8892
    "M136" = -12831.8892
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$clears = @("N45")
foreach ($cellRef in $clears) {
    $ws.Range($cellRef).ClearContents()
}

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$updates = @{
    "H134" = 4593.5
    "I134" = 4887.409
    "J134" = 2977
    "K134" = 14662.227
    "L134" = 8931
    "M134" = -12127.227
    "N134" = -14001
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$updates = @{
    "H31" = 4397.0835
    "I31" = 2524.0908
    "K31" = 2524.0908
    "M31" = -2229.0908
    "H32" = 4119.8
    "I32" = 4119.8
    "J32" = 0
    "K32" = 4119.8
    "L32" = 0
    "N32" = -3803.8
    "H34" = 4397.0835
    "I34" = 2524.0908
    "K34" = 2524.0908
    "M34" = -2322.0908
    "H105" = 1038.4286
    "I105" = 963.8
    "K105" = 963.8
    "M105" = 783.2
    "H132" = 1719.7368
    "I132" = 1716.9286
    "K132" = 5150.7858
    "M132" = -2620.7858
    "H134" = 1890.45
    "I134" = 1852
    "K134" = 5556
    "M134" = -3021
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$clears = @("M32")
foreach ($cellRef in $clears) {
    $ws.Range($cellRef).ClearContents()
}

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$updates = @{
    "H12" = 110.333336
    "I12" = 125.5
    "K12" = 376.5
    "M12" = -203.5
    "H63" = 5600
    "I63" = 1200
    "J63" = 10000
    "K63" = 3600
    "L63" = 30000
    "M63" = -2851
    "N63" = -31498
    "H66" = 5600
    "I66" = 1200
    "J66" = 10000
    "K66" = 10800
    "L66" = 90000
    "M66" = -7056
    "N66" = -97488
    "I107" = 140
    "J107" = 150
    "K107" = 420
    "L107" = 450
    "M107" = 1500
    "N107" = -4290
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$updates = @{
    "H43" = 9979.25
    "I43" = 2465
    "J43" = 22503
    "K43" = 2465
    "L43" = 22503
    "M43" = -2314
    "N43" = -22805
    "H80" = 2479.4
    "I80" = 2849.25
    "J80" = 1000
    "K80" = 2849.25
    "L80" = 1000
    "M80" = -1851.25
    "N80" = -2996
    "H83" = 2479.4
    "I83" = 2849.25
    "J83" = 1000
    "K83" = 14246.25
    "L83" = 5000
    "M83" = -9254.25
    "N83" = -14984
    "H132" = 1905.5454
    "I132" = 1796.1
    "J132" = 3000
    "K132" = 5388.299999999999
    "L132" = 9000
    "M132" = -2858.299999999999
    "N132" = -14060
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$updates = @{
    "H22" = 696.6667
    "I22" = 445
    "J22" = 1200
    "K22" = 445
    "L22" = 1200
    "M22" = -150
    "N22" = -1790
    "H27" = 696.6667
    "I27" = 445
    "J27" = 1200
    "K27" = 445
    "L27" = 1200
    "M27" = -338
    "N27" = -1414
    "H46" = 1941
    "I46" = 2271.5
    "J46" = 949.5
    "K46" = 2271.5
    "L46" = 949.5
    "M46" = -2083.5
    "N46" = -1325.5
    "H68" = 2490.6667
    "I68" = 2398.5
    "K68" = 2398.5
    "M68" = -1649.5
    "H71" = 2490.6667
    "I71" = 2398.5
    "K71" = 11992.5
    "M71" = -8248.5
    "H132" = 2517.3157
    "I132" = 2548.9412
    "J132" = 2248.5
    "K132" = 7646.823600000001
    "L132" = 6745.5
    "M132" = -5116.823600000001
    "N132" = -11805.5
    "H136" = 9110.333000000001
    "I136" = 6213.2856
    "K136" = 18639.8568
    "M136" = -16089.8568
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$updates = @{
    "H20" = 0
    "I20" = 0
    "K20" = 0
    "H45" = 48152
    "I45" = 40008
    "J45" = 54938.668
    "K45" = 40008
    "L45" = 54938.668
    "M45" = -39517
    "N45" = -55920.668
    "H81" = 2925.9092
    "I81" = 2925.9092
    "K81" = 5851.8184
    "M81" = -4790.8184
    "H84" = 2925.9092
    "I84" = 2925.9092
    "K84" = 29259.092
    "M84" = -23955.092
    "H113" = 547.3125
    "I113" = 413.16666
    "K113" = 1239.49998
    "M113" = 930.5000199999999
    "H126" = 1752
    "I126" = 1004
    "J126" = 2500
    "K126" = 3012
    "L126" = 7500
    "M126" = -542
    "N126" = -12440
    "H132" = 2555.7354
    "I132" = 2330.1516
    "K132" = 6990.4548
    "M132" = -4460.4548
    "H135" = 200715
    "J135" = 200715
    "L135" = 200715
    "N135" = -210855
    "H136" = 28314.834
    "I136" = 27998
    "J136" = 28948.5
    "K136" = 83994
    "L136" = 86845.5
    "M136" = -81444
    "N136" = -91945.5
    "H141" = 54997.5
    "J141" = 54997.5
    "L141" = 54997.5
    "N141" = -65357.5
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$clears = @("M20")
foreach ($cellRef in $clears) {
    $ws.Range($cellRef).ClearContents()
}
